$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7085713744163513
$ws.Range("B1").Value = 3.703454971313477
$ws.Range("C1").Value = 6.710021495819092
$ws.Range("D1").Value = 5.046497344970703
$ws.Range("E1").Value = 2.393457651138306
